$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.871.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.71%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.607.89'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.92%  '

$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.38%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.612.05'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.09%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("E9").Value = '  +1.08%  '

$ws.Range("E10").Value = '  +0.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.27'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.394'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.217.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.23%  '

$ws.Range("E15").Value = '  +0.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.602.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.47%  '

$ws.Range("E17").Value = '  +0.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.948.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.67%  '

$ws.Range("E21").Value = '  +0.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '399.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.83%  '

$ws.Range("E23").Value = '  +3.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.749.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("E27").Value = '  +1.67%  '

$ws.Range("E28").Value = '  +3.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.68'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +29.95%  '

$ws.Range("E30").Value = '  +4.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.72'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.68%  '

$ws.Range("E32").Value = '  +0.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.602.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.70'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.92%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.149'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.83%  '

$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.54%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.63%  '

$ws.Range("E39").Value = '  +2.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '171.97'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0836'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.69%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.845'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.63%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.44'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.55'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.56%  '

$ws.Range("E47").Value = '  -0.23%  '

$ws.Range("E48").Value = '  +1.17%  '

$ws.Range("E49").Value = '  +4.00%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.470.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.41%  '

$ws.Range("E51").Value = '  +3.22%  '
